$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCAUY")

$ws.Range("D8").Value = 787300
$ws.Range("E8").Value = 760600
$ws.Range("F8").Value = 721700
$ws.Range("G8").Value = 818500
$ws.Range("H8").Value = 905800
$ws.Range("I8").Value = 878000
$ws.Range("J8").Value = 956200
$ws.Range("D9").Value = 759600
$ws.Range("E9").Value = 735200
$ws.Range("F9").Value = 690800
$ws.Range("G9").Value = 735000
$ws.Range("H9").Value = 803900
$ws.Range("I9").Value = 774700
$ws.Range("J9").Value = 829200
$ws.Range("D10").Value = 27600
$ws.Range("E10").Value = 25500
$ws.Range("F10").Value = 30900
$ws.Range("G10").Value = 83500
$ws.Range("H10").Value = 101900
$ws.Range("I10").Value = 103300
$ws.Range("J10").Value = 127000
$ws.Range("H12").Value = 600
$ws.Range("D14").Value = 105300
$ws.Range("I14").Value = 3100
$ws.Range("J14").Value = 3400
$ws.Range("D17").Value = 996900
$ws.Range("E17").Value = 862800
$ws.Range("F17").Value = 806500
$ws.Range("G17").Value = 858000
$ws.Range("H17").Value = 932300
$ws.Range("I17").Value = 537400
$ws.Range("J17").Value = 666900
$ws.Range("D18").Value = -209700
$ws.Range("E18").Value = -102100
$ws.Range("F18").Value = -84800
$ws.Range("G18").Value = -39500
$ws.Range("H18").Value = -26400
$ws.Range("I18").Value = 340500
$ws.Range("J18").Value = 289300
$ws.Range("D20").Value = 808900
$ws.Range("E20").Value = 630100
$ws.Range("F20").Value = 600000
$ws.Range("G20").Value = 855700
$ws.Range("H20").Value = 540400
$ws.Range("D21").Value = 640200
$ws.Range("E21").Value = 567800
$ws.Range("F21").Value = 542000
$ws.Range("G21").Value = 836600
$ws.Range("H21").Value = 534500
$ws.Range("I21").Value = 361200
$ws.Range("J21").Value = "NA"
$ws.Range("D22").Value = 20500
$ws.Range("E22").Value = 19800
$ws.Range("F22").Value = 21800
$ws.Range("G22").Value = 23200
$ws.Range("H22").Value = 20600
$ws.Range("D23").Value = 578800
$ws.Range("E23").Value = 508200
$ws.Range("F23").Value = 493400
$ws.Range("G23").Value = 792900
$ws.Range("H23").Value = 493400
$ws.Range("I23").Value = 340500
$ws.Range("J23").Value = 289300
$ws.Range("D24").Value = 5000
$ws.Range("E24").Value = 5300
$ws.Range("F24").Value = 6600
$ws.Range("G24").Value = 6400
$ws.Range("I24").Value = 8500
$ws.Range("J24").Value = 8600
$ws.Range("D26").Value = 573700
$ws.Range("E26").Value = 502900
$ws.Range("F26").Value = 486800
$ws.Range("G26").Value = 786600
$ws.Range("H26").Value = 492200
$ws.Range("I26").Value = 332000
$ws.Range("J26").Value = 280700
$ws.Range("D27").Value = 649500
$ws.Range("E27").Value = 546500
$ws.Range("F27").Value = 518700
$ws.Range("G27").Value = 801900
$ws.Range("H27").Value = 500800
$ws.Range("I27").Value = 341500
$ws.Range("J27").Value = 269000
$ws.Range("D32").Value = -808900
$ws.Range("E32").Value = -630100
$ws.Range("F32").Value = -600000
$ws.Range("G32").Value = -855700
$ws.Range("H32").Value = -540400
$ws.Range("D33").Value = 649500
$ws.Range("E33").Value = 546500
$ws.Range("F33").Value = 518700
$ws.Range("G33").Value = 801900
$ws.Range("H33").Value = 500800
$ws.Range("I33").Value = 341500
$ws.Range("J33").Value = 269000
$ws.Range("D35").Value = 649500
$ws.Range("E35").Value = 546500
$ws.Range("F35").Value = 518700
$ws.Range("G35").Value = 801900
$ws.Range("H35").Value = 500800
$ws.Range("I35").Value = 341500
$ws.Range("J35").Value = 269000
$ws.Range("D41").Value = 266300
$ws.Range("E41").Value = 139600
$ws.Range("F41").Value = 158900
$ws.Range("G41").Value = 174900
$ws.Range("I41").Value = 124100
$ws.Range("J41").Value = 86900
$ws.Range("D42").Value = 260800
$ws.Range("E42").Value = 227400
$ws.Range("F42").Value = 297000
$ws.Range("G42").Value = 199900
$ws.Range("I42").Value = 165100
$ws.Range("J42").Value = 190400
$ws.Range("D43").Value = 644900
$ws.Range("E43").Value = 501600
$ws.Range("F43").Value = 419700
$ws.Range("G43").Value = 436600
$ws.Range("I43").Value = 343400
$ws.Range("J43").Value = 293100
$ws.Range("D44").Value = 154900
$ws.Range("E44").Value = 163900
$ws.Range("F44").Value = 179700
$ws.Range("G44").Value = 118200
$ws.Range("I44").Value = 124400
$ws.Range("J44").Value = 109400
$ws.Range("D45").Value = 13600
$ws.Range("E45").Value = 7800
$ws.Range("F45").Value = 9500
$ws.Range("G45").Value = 11900
$ws.Range("I45").Value = 195300
$ws.Range("J45").Value = 215300
$ws.Range("D46").Value = 1340400
$ws.Range("E46").Value = 1040300
$ws.Range("F46").Value = 1064800
$ws.Range("G46").Value = 941600
$ws.Range("I46").Value = 952400
$ws.Range("J46").Value = 895200
$ws.Range("D47").Value = 3682400
$ws.Range("E47").Value = 2930000
$ws.Range("F47").Value = 2334800
$ws.Range("G47").Value = 1899100
$ws.Range("I47").Value = 1107100
$ws.Range("J47").Value = 710300
$ws.Range("D48").Value = 381000
$ws.Range("E48").Value = 333900
$ws.Range("F48").Value = 303100
$ws.Range("G48").Value = 290900
$ws.Range("I48").Value = 259000
$ws.Range("J48").Value = 247800
$ws.Range("D49").Value = 103300
$ws.Range("E49").Value = 198700
$ws.Range("F49").Value = 211200
$ws.Range("G49").Value = 147800
$ws.Range("I49").Value = 62900
$ws.Range("J49").Value = 38800
$ws.Range("D52").Value = 111100
$ws.Range("E52").Value = 104200
$ws.Range("F52").Value = 103900
$ws.Range("G52").Value = 164800
$ws.Range("I52").Value = 1700
$ws.Range("J52").Value = 9200
$ws.Range("D54").Value = 5618200
$ws.Range("E54").Value = 4607100
$ws.Range("F54").Value = 4017800
$ws.Range("G54").Value = 3444100
$ws.Range("I54").Value = 2383100
$ws.Range("J54").Value = 1901200
$ws.Range("D57").Value = 486600
$ws.Range("E57").Value = 493300
$ws.Range("F57").Value = 450900
$ws.Range("G57").Value = 439800
$ws.Range("I57").Value = 463000
$ws.Range("J57").Value = 365900
$ws.Range("D58").Value = 829700
$ws.Range("E58").Value = 542400
$ws.Range("F58").Value = 555400
$ws.Range("G58").Value = 478300
$ws.Range("I58").Value = 419600
$ws.Range("J58").Value = 458200
$ws.Range("D59").Value = 311000
$ws.Range("E59").Value = 199400
$ws.Range("F59").Value = 162000
$ws.Range("G59").Value = 140600
$ws.Range("I59").Value = 135100
$ws.Range("J59").Value = 151300
$ws.Range("D60").Value = 1627300
$ws.Range("E60").Value = 1235200
$ws.Range("F60").Value = 1168300
$ws.Range("G60").Value = 1058800
$ws.Range("I60").Value = 1017700
$ws.Range("J60").Value = 975300
$ws.Range("D61").Value = 11900
$ws.Range("D62").Value = 16500
$ws.Range("E62").Value = 18100
$ws.Range("F62").Value = 20300
$ws.Range("G62").Value = 17700
$ws.Range("D66").Value = 1682000
$ws.Range("E66").Value = 1086200
$ws.Range("F66").Value = 1065100
$ws.Range("G66").Value = 931400
$ws.Range("I66").Value = 896800
$ws.Range("J66").Value = 863900
$ws.Range("D72").Value = 3794600
$ws.Range("E72").Value = 3215500
$ws.Range("F72").Value = 2739900
$ws.Range("G72").Value = 2288400
$ws.Range("I72").Value = 1102300
$ws.Range("J72").Value = 654900
$ws.Range("D76").Value = 3936300
$ws.Range("E76").Value = 3520900
$ws.Range("F76").Value = 2952700
$ws.Range("G76").Value = 2512800
$ws.Range("I76").Value = 1486300
$ws.Range("J76").Value = 1037300
$ws.Range("D81").Value = 649500
$ws.Range("E81").Value = 546500
$ws.Range("F81").Value = 518700
$ws.Range("G81").Value = 801900
$ws.Range("H81").Value = 500800
$ws.Range("I81").Value = 341500
$ws.Range("J81").Value = 269000
$ws.Range("D83").Value = 40900
$ws.Range("E83").Value = 39700
$ws.Range("F83").Value = 26700
$ws.Range("G83").Value = 20500
$ws.Range("H83").Value = 20400
$ws.Range("I83").Value = 20600
$ws.Range("J83").Value = "NA"
$ws.Range("D89").Value = -378000
$ws.Range("E89").Value = -128300
$ws.Range("F89").Value = -151500
$ws.Range("G89").Value = 114200
$ws.Range("H89").Value = -22800
$ws.Range("J89").Value = -93100
$ws.Range("D91").Value = -92600
$ws.Range("E91").Value = -71300
$ws.Range("F91").Value = -90000
$ws.Range("G91").Value = -109400
$ws.Range("H91").Value = -91400
$ws.Range("I91").Value = -78100
$ws.Range("J91").Value = -45900
$ws.Range("D94").Value = 189700
$ws.Range("E94").Value = 169200
$ws.Range("F94").Value = 53500
$ws.Range("G94").Value = 31800
$ws.Range("H94").Value = 69800
$ws.Range("I94").Value = 53800
$ws.Range("J94").Value = "NA"
$ws.Range("D96").Value = -41100
$ws.Range("E96").Value = -70500
$ws.Range("F96").Value = -67200
$ws.Range("G96").Value = -64900
$ws.Range("H96").Value = -58400
$ws.Range("D100").Value = 306300
$ws.Range("E100").Value = -60700
$ws.Range("F100").Value = 82000
$ws.Range("G100").Value = -105200
$ws.Range("H100").Value = -37000
$ws.Range("I100").Value = -16400
$ws.Range("J100").Value = "NA"
$ws.Range("J101").Value = "NA"
$ws.Range("D102").Value = 118000
$ws.Range("E102").Value = -19900
$ws.Range("F102").Value = -16000
$ws.Range("G102").Value = 40900
$ws.Range("H102").Value = 9900
$ws.Range("I102").Value = 37200
$ws.Range("J102").Value = 23400
